# Update the division-problem values in the worksheet table.
# Cells are addressed by (row, column) to avoid ambiguity since several
# "find" strings (e.g. "36÷3=", "81÷9=") occur more than once in the
# document but must map to different replacement values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "81÷8=" },
    @{ Row = 1;  Col = 2; New = "51÷4=" },
    @{ Row = 1;  Col = 3; New = "78÷2=" },
    @{ Row = 1;  Col = 4; New = "24÷7=" },
    @{ Row = 1;  Col = 5; New = "39÷5=" },

    @{ Row = 5;  Col = 1; New = "57÷8=" },
    @{ Row = 5;  Col = 2; New = "43÷8=" },
    @{ Row = 5;  Col = 3; New = "58÷5=" },
    @{ Row = 5;  Col = 4; New = "35÷5=" },
    @{ Row = 5;  Col = 5; New = "71÷4=" },

    @{ Row = 9;  Col = 1; New = "80÷4=" },
    @{ Row = 9;  Col = 2; New = "51÷7=" },
    @{ Row = 9;  Col = 3; New = "81÷4=" },
    @{ Row = 9;  Col = 4; New = "83÷8=" },
    @{ Row = 9;  Col = 5; New = "79÷6=" },

    @{ Row = 13; Col = 1; New = "39÷9=" },
    @{ Row = 13; Col = 2; New = "69÷6=" },
    @{ Row = 13; Col = 3; New = "88÷8=" },
    @{ Row = 13; Col = 4; New = "67÷7=" },
    @{ Row = 13; Col = 5; New = "36÷2=" },

    @{ Row = 17; Col = 1; New = "39÷6=" },
    @{ Row = 17; Col = 2; New = "29÷5=" },
    @{ Row = 17; Col = 3; New = "65÷6=" },
    @{ Row = 17; Col = 4; New = "76÷4=" },
    @{ Row = 17; Col = 5; New = "15÷5=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
